$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$m.ApplyTheme("theme1.xml")
